$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"
$ws.Range("B2").Style = "Normal"
$ws.Range("D2").Value = 0.31435
$ws.Range("E2").Value = 0.1752
$ws.Range("G2").Value = 0.08872023030620257
$ws.Range("H2").Value = 0.08872023030620257
$ws.Range("I2").Value = 0.09604815493326356
$ws.Range("J2").Value = 0.08401732774448746
$ws.Range("K2").Value = 29.09
$ws.Range("L2").Value = 0.07613190264328711
$ws.Range("M2").Value = 6.179079999999999
$ws.Range("N2").Value = 0.03751718275652702
$ws.Range("O2").Value = 0.2124125128910278
$ws.Range("P2").Value = 6.179079999999999
$ws.Range("Q2").Value = 0.03751718275652702
$ws.Range("R2").Value = 0.2124125128910278
$ws.Range("U2").Value = 121.3
$ws.Range("V2").Value = 0.7364905889496055
$ws.Range("W2").Value = 0.1657556461198845
$ws.Range("X2").Value = 0.05509253993312482
$ws.Range("Y2").Value = 0.1106631061867597
$ws.Range("Z2").Value = 2.934715821812596
$ws.Range("AA2").Value = 0.2426203787389128
$ws.Range("AB2").Value = 0.05199192126089672
$ws.Range("AC2").Value = 0.190628457478016
$ws.Range("AD2").Value = 25.4
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 25.4
$ws.Range("AG2").Value = -95.90000000000001
$ws.Range("AH2").Value = 0.1336138874276696
$ws.Range("AI2").Value = 0.1157702825888788
$ws.Range("AJ2").Value = -1.39389534883721
$ws.Range("AK2").Value = -0.9775739041794089
$ws.Range("AL2").Value = 0.7769999999999999
$ws.Range("AM2").Value = 0.7769999999999999
$ws.Range("AN2").Value = 0.6446700507614213
$ws.Range("AO2").Value = 47.23294723294724
$ws.Range("AP2").Value = -2.434010152284264
$ws.Range("AQ2").Value = 47.23294723294724

# --- Row 3 updates ---
$ws.Range("B3").Value = "Swan General Ltd (MUSE:SWAN.N0000)"
$ws.Range("D3").Value = 0.529
$ws.Range("E3").Value = 0.291
$ws.Range("G3").Value = 0.07052510540436949
$ws.Range("H3").Value = 0.07052510540436949
$ws.Range("I3").Value = 0.09237255653507093
$ws.Range("J3").Value = 0.08579505542770985
$ws.Range("K3").Value = 21
$ws.Range("L3").Value = 0.08049060942889996
$ws.Range("M3").Value = 2.98908
$ws.Range("N3").Value = 0.04782527999999999
$ws.Range("O3").Value = 0.1423371428571428
$ws.Range("P3").Value = 2.98908
$ws.Range("Q3").Value = 0.04782527999999999
$ws.Range("R3").Value = 0.1423371428571428
$ws.Range("U3").Value = 98.90000000000001
$ws.Range("V3").Value = 1.5824
$ws.Range("W3").Value = 0.2243589743589744
$ws.Range("X3").Value = 0.05071774044130588
$ws.Range("Y3").Value = 0.1736412339176685
$ws.Range("Z3").Value = 3.310913705583756
$ws.Range("AA3").Value = 0.2840600248869226
$ws.Range("AB3").Value = 0.05071774044130588
$ws.Range("AC3").Value = 0.2333422844456167
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -98.90000000000001
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 2.717032967032967
$ws.Range("AK3").Value = -30.90625000000011
$ws.Range("AL3").Value = 0.061
$ws.Range("AM3").Value = 0.061
$ws.Range("AO3").Value = 395.0819672131148
$ws.Range("AP3").Value = -3.893700787401575
$ws.Range("AQ3").Value = 395.0819672131148

# --- Row 4 updates ---
$ws.Range("B4").Value = "MUA Ltd (MUSE:MUAL.N0000)"
$ws.Range("D4").Value = 0.09970000000000001
$ws.Range("E4").Value = 0.0594
$ws.Range("G4").Value = 0.1278877887788779
$ws.Range("H4").Value = 0.1278877887788779
$ws.Range("I4").Value = 0.103960396039604
$ws.Range("J4").Value = 0.08531922157733014
$ws.Range("K4").Value = 8.09
$ws.Range("L4").Value = 0.06674917491749174
$ws.Range("M4").Value = 3.19
$ws.Range("N4").Value = 0.0312133072407045
$ws.Range("O4").Value = 0.3943139678615575
$ws.Range("P4").Value = 3.19
$ws.Range("Q4").Value = 0.0312133072407045
$ws.Range("R4").Value = 0.3943139678615575
$ws.Range("U4").Value = 22.4
$ws.Range("V4").Value = 0.2191780821917808
$ws.Range("W4").Value = 0.1071523178807947
$ws.Range("X4").Value = 0.05946733942494377
$ws.Range("Y4").Value = 0.04768497845585093
$ws.Range("Z4").Value = 2.357976653696498
$ws.Range("AA4").Value = 0.201180732590903
$ws.Range("AB4").Value = 0.05326610208048757
$ws.Range("AC4").Value = 0.1479146305104154
$ws.Range("AD4").Value = 25.4
$ws.Range("AF4").Value = 25.4
$ws.Range("AG4").Value = 3
$ws.Range("AH4").Value = 0.1990595611285266
$ws.Range("AI4").Value = 0.216538789428815
$ws.Range("AJ4").Value = 0.0285171102661597
$ws.Range("AK4").Value = 0.03161222339304531
$ws.Range("AL4").Value = 0.716
$ws.Range("AM4").Value = 0.716
$ws.Range("AN4").Value = 1.814285714285714
$ws.Range("AO4").Value = 17.59776536312849
$ws.Range("AP4").Value = 0.2142857142857143
$ws.Range("AQ4").Value = 17.59776536312849

# --- Remove row 5 (Eagle Insurance Limited) ---
$ws.Rows(5).Delete()
